{"js": "// Replace the date line and all the three-digit x one-digit multiplication\n// answers in the table with their updated values, per the commit diff.\nconst replacements = [\n  [\"2024-12-19 Thursday\", \"2024-12-20 Friday\"],\n  [\"391\u00d75=1955\", \"491\u00d75=2455\"],\n  [\"207\u00d79=1863\", \"721\u00d78=5768\"],\n  [\"690\u00d73=2070\", \"304\u00d76=1824\"],\n  [\"737\u00d74=2948\", \"532\u00d77=3724\"],\n  [\"596\u00d78=4768\", \"185\u00d78=1480\"],\n  [\"104\u00d74=416\", \"474\u00d77=3318\"],\n  [\"622\u00d77=4354\", \"640\u00d75=3200\"],\n  [\"640\u00d79=5760\", \"174\u00d73=522\"],\n  [\"174\u00d75=870\", \"762\u00d76=4572\"],\n  [\"705\u00d73=2115\", \"178\u00d77=1246\"],\n  [\"483\u00d76=2898\", \"642\u00d72=1284\"],\n  [\"196\u00d77=1372\", \"844\u00d79=7596\"],\n  [\"952\u00d76=5712\", \"845\u00d74=3380\"],\n  [\"213\u00d75=1065\", \"884\u00d73=2652\"],\n  [\"580\u00d78=4640\", \"986\u00d79=8874\"],\n  [\"914\u00d73=2742\", \"507\u00d79=4563\"],\n  [\"510\u00d77=3570\", \"210\u00d77=1470\"],\n  [\"101\u00d78=808\", \"329\u00d76=1974\"],\n  [\"424\u00d77=2968\", \"780\u00d79=7020\"],\n  [\"796\u00d72=1592\", \"330\u00d77=2310\"],\n  [\"245\u00d79=2205\", \"931\u00d73=2793\"],\n  [\"934\u00d79=8406\", \"599\u00d72=1198\"],\n  [\"697\u00d77=4879\", \"612\u00d73=1836\"],\n  [\"668\u00d76=4008\", \"551\u00d74=2204\"],\n  [\"782\u00d78=6256\", \"643\u00d72=1286\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all the three-digit x one-digit multiplication\n# answers in the table with their updated values, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-19 Thursday\", \"2024-12-20 Friday\"),\n    @(\"391\u00d75=1955\", \"491\u00d75=2455\"),\n    @(\"207\u00d79=1863\", \"721\u00d78=5768\"),\n    @(\"690\u00d73=2070\", \"304\u00d76=1824\"),\n    @(\"737\u00d74=2948\", \"532\u00d77=3724\"),\n    @(\"596\u00d78=4768\", \"185\u00d78=1480\"),\n    @(\"104\u00d74=416\", \"474\u00d77=3318\"),\n    @(\"622\u00d77=4354\", \"640\u00d75=3200\"),\n    @(\"640\u00d79=5760\", \"174\u00d73=522\"),\n    @(\"174\u00d75=870\", \"762\u00d76=4572\"),\n    @(\"705\u00d73=2115\", \"178\u00d77=1246\"),\n    @(\"483\u00d76=2898\", \"642\u00d72=1284\"),\n    @(\"196\u00d77=1372\", \"844\u00d79=7596\"),\n    @(\"952\u00d76=5712\", \"845\u00d74=3380\"),\n    @(\"213\u00d75=1065\", \"884\u00d73=2652\"),\n    @(\"580\u00d78=4640\", \"986\u00d79=8874\"),\n    @(\"914\u00d73=2742\", \"507\u00d79=4563\"),\n    @(\"510\u00d77=3570\", \"210\u00d77=1470\"),\n    @(\"101\u00d78=808\", \"329\u00d76=1974\"),\n    @(\"424\u00d77=2968\", \"780\u00d79=7020\"),\n    @(\"796\u00d72=1592\", \"330\u00d77=2310\"),\n    @(\"245\u00d79=2205\", \"931\u00d73=2793\"),\n    @(\"934\u00d79=8406\", \"599\u00d72=1198\"),\n    @(\"697\u00d77=4879\", \"612\u00d73=1836\"),\n    @(\"668\u00d76=4008\", \"551\u00d74=2204\"),\n    @(\"782\u00d78=6256\", \"643\u00d72=1286\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
